$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to start with a merged-looking title row ("Bảng Tra Tiết
# Diện Dây Dẫn Tham Khảo ...") in row 1, followed by the real header row
# (Tiết diện / Khả năng chịu tải ...) in row 2 and the data rows below it.
# The edit removes that title row entirely, so the header row becomes row 1
# and the data rows shift up to rows 2-17.
$ws.Rows("1:1").Delete() | Out-Null

# The saved workbook's selection/active cell is A6 (not the old C19).
$ws.Range("A6").Select() | Out-Null
